$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group Member")

# --- Activate the "Group Member" sheet (was "Daily level") ---
$ws.Activate()

# --- Column A width: 15.26953125 -> 11 (stored width units) ---
$ws.Columns.Item(1).ColumnWidth = 10.166666666666666

# --- Target state for columns B, C, D on rows 2..11 ---
# type: "s" = string value, "d" = date serial number value
# every one of these cells ends up either:
#   - General number format + centered horizontal alignment (text cells), or
#   - m/d/yyyy date format + centered horizontal alignment (date cells)

$rows = @(
    @{ r=2;  B=@{t="s"; v="Pending"};  C=@{t="s"; v="Pending"};  D=@{t="s"; v="pending"} },
    @{ r=3;  B=@{t="s"; v="Pending"};  C=@{t="s"; v="Pending"};  D=@{t="s"; v="pending"} },
    @{ r=4;  B=@{t="d"; v=44064};      C=@{t="s"; v="Accepted"}; D=@{t="s"; v="Accepted"} },
    @{ r=5;  B=@{t="d"; v=44044};      C=@{t="s"; v="Accepted"}; D=@{t="s"; v="Accepted"} },
    @{ r=6;  B=@{t="d"; v=44065};      C=@{t="s"; v="Accepted"}; D=@{t="s"; v="Accepted"} },
    @{ r=7;  B=@{t="s"; v="Pending"};  C=@{t="s"; v="Pending"};  D=@{t="s"; v="pending"} },
    @{ r=8;  B=@{t="d"; v=44044};      C=@{t="s"; v="Accepted"}; D=@{t="d"; v=44063} },
    @{ r=9;  B=@{t="d"; v=44048};      C=@{t="s"; v="Accepted"}; D=@{t="s"; v="Accepted"} },
    @{ r=10; B=@{t="d"; v=44063};      C=@{t="s"; v="Accepted"}; D=@{t="s"; v="Accepted"} },
    @{ r=11; B=@{t="s"; v="Rejected"}; C=@{t="d"; v=44064};      D=@{t="s"; v="Accepted"} }
)

# Build the two required styles exactly once each (avoids creating stray,
# unused intermediate cellXfs entries), then reuse them via PasteSpecial.
$dateStyleSourceSet = $false
$dateStyleSource = $null

foreach ($row in $rows) {
    $r = $row.r
    foreach ($col in @("B","C","D")) {
        $cellInfo = $row[$col]
        $cell = $ws.Range("$col$r")
        $cell.Value = $cellInfo.v
        if ($cellInfo.t -eq "d") {
            if (-not $dateStyleSourceSet) {
                $cell.NumberFormat = "m/d/yyyy"
                $cell.HorizontalAlignment = -4108
                $dateStyleSource = $cell
                $dateStyleSourceSet = $true
            } else {
                $dateStyleSource.Copy()
                $cell.PasteSpecial(-4122)
            }
        } else {
            $cell.HorizontalAlignment = -4108
        }
    }
}

# --- Update selection on the Group Member sheet ---
$ws.Range("E2").Select()
